# Horarios actualizados Línea 141 - 820
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with a
# fresh scrape batch (timestamp 05:21:16): new rows are inserted/appended and
# the "Minutos" countdown of previously-scraped rows that are still pending is
# refreshed against the new scrape time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 05:21:16"
$ws1.Range("A3").Value = "Total filas: 44"

# Insert a brand-new scraped row at position 23 (pushes the existing rows
# 23-43 down to 24-44).
$ws1.Rows.Item(23).Insert()
$ws1.Cells.Item(23, 1).Value = "05:21:16"
$ws1.Cells.Item(23, 2).Value = "05:26"
$ws1.Cells.Item(23, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(23, 4).Value = 5
$ws1.Cells.Item(23, 5).Value = "LP1912"

# Refresh Hora_Scrap (A) and Minutos (D) for the rows that were re-scraped
# at 05:21:16 (their positions shifted down by one after the insert above).
$ws1.Cells.Item(26, 1).Value = "05:21:16"
$ws1.Cells.Item(26, 4).Value = 14

$ws1.Cells.Item(28, 1).Value = "05:21:16"
$ws1.Cells.Item(28, 4).Value = 25

$ws1.Cells.Item(29, 1).Value = "05:21:16"
$ws1.Cells.Item(29, 4).Value = 33

$ws1.Cells.Item(30, 1).Value = "05:21:16"
$ws1.Cells.Item(30, 4).Value = 43

$ws1.Cells.Item(31, 1).Value = "05:21:16"
$ws1.Cells.Item(31, 4).Value = 50

$ws1.Cells.Item(33, 1).Value = "05:21:16"
$ws1.Cells.Item(33, 4).Value = 53

$ws1.Cells.Item(35, 1).Value = "05:21:16"
$ws1.Cells.Item(35, 4).Value = 60

$ws1.Cells.Item(37, 1).Value = "05:21:16"
$ws1.Cells.Item(37, 4).Value = 66

$ws1.Cells.Item(39, 1).Value = "05:21:16"
$ws1.Cells.Item(39, 4).Value = 69

$ws1.Cells.Item(40, 1).Value = "05:21:16"
$ws1.Cells.Item(40, 4).Value = 70

$ws1.Cells.Item(42, 1).Value = "05:21:16"
$ws1.Cells.Item(42, 4).Value = 83

$ws1.Cells.Item(43, 1).Value = "05:21:16"
$ws1.Cells.Item(43, 4).Value = 85

# Append the newly-scraped rows at the bottom of the list (45-49).
$ws1.Cells.Item(45, 1).Value = "05:21:16"
$ws1.Cells.Item(45, 2).Value = "07:00"
$ws1.Cells.Item(45, 3).Value = "14_ABASTO"
$ws1.Cells.Item(45, 4).Value = 99
$ws1.Cells.Item(45, 5).Value = "LP1912"

$ws1.Cells.Item(46, 1).Value = "05:21:16"
$ws1.Cells.Item(46, 2).Value = "07:05"
$ws1.Cells.Item(46, 3).Value = "15_ABASTO"
$ws1.Cells.Item(46, 4).Value = 104
$ws1.Cells.Item(46, 5).Value = "LP1912"

$ws1.Cells.Item(47, 1).Value = "05:21:16"
$ws1.Cells.Item(47, 2).Value = "07:07"
$ws1.Cells.Item(47, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(47, 4).Value = 106
$ws1.Cells.Item(47, 5).Value = "LP1912"

$ws1.Cells.Item(48, 1).Value = "05:21:16"
$ws1.Cells.Item(48, 2).Value = "07:11"
$ws1.Cells.Item(48, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(48, 4).Value = 110
$ws1.Cells.Item(48, 5).Value = "LP1912"

$ws1.Cells.Item(49, 1).Value = "05:21:16"
$ws1.Cells.Item(49, 2).Value = "07:16"
$ws1.Cells.Item(49, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(49, 4).Value = 115
$ws1.Cells.Item(49, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:21:16"
$ws2.Range("A3").Value = "Total filas: 12"

$ws2.Cells.Item(13, 1).Value = "05:21:16"
$ws2.Cells.Item(13, 4).Value = 14

$ws2.Cells.Item(14, 1).Value = "05:21:16"
$ws2.Cells.Item(14, 4).Value = 50

$ws2.Cells.Item(15, 1).Value = "05:21:16"
$ws2.Cells.Item(15, 4).Value = 85

# Append new row 17 with the latest scrape.
$ws2.Cells.Item(17, 1).Value = "05:21:16"
$ws2.Cells.Item(17, 2).Value = "07:11"
$ws2.Cells.Item(17, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(17, 4).Value = 110
$ws2.Cells.Item(17, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 05:21:16"
$ws3.Range("A3").Value = "Total filas: 8"

$ws3.Cells.Item(8, 1).Value = "05:21:16"
$ws3.Cells.Item(8, 4).Value = 23

$ws3.Cells.Item(10, 1).Value = "05:21:16"
$ws3.Cells.Item(10, 4).Value = 48

$ws3.Cells.Item(12, 1).Value = "05:21:16"
$ws3.Cells.Item(12, 4).Value = 72

# Append new row 13 with the latest scrape.
$ws3.Cells.Item(13, 1).Value = "05:21:16"
$ws3.Cells.Item(13, 2).Value = "07:00"
$ws3.Cells.Item(13, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(13, 4).Value = 99
$ws3.Cells.Item(13, 5).Value = "L6173"

Write-Output "Horarios actualizados - Linea 141 - 820"
